# Generic Utility Implementation contd...
#
# Update two test-data cells on the "ContactsTestData" sheet and switch the
# workbook's active sheet/selection/zoom to match the author's last-saved
# view state.

$wb = $excel.ActiveWorkbook

$contactsSheet = $wb.Worksheets.Item("ContactsTestData")

# --- Content edits on ContactsTestData ---------------------------------
# Row 6 "Company" test value: Microsoft -> ABC
$contactsSheet.Range("D6").Value = "ABC"
# Row 5 "Last Name" test value: Abc -> XYZ
$contactsSheet.Range("D5").Value = "XYZ"

# --- View-state edits ----------------------------------------------------
# Make ContactsTestData the active/selected tab (was OrganizationsTestData).
$contactsSheet.Activate()

# Restore the zoom level and selection the author left the sheet at.
$excel.ActiveWindow.Zoom = 183
$contactsSheet.Range("E5").Select()
